$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Accounts")

# New Jodine accounts being added to the budget
$ws.Range("A16").Value = "Jodine Transactional"
$ws.Range("C16").Value = "Transactional"
$ws.Range("D16").Value = "Jodine"

$ws.Range("A17").Value = "Jodine Cash"
$ws.Range("C17").Value = "Transactional"
$ws.Range("D17").Value = "Jodine"

$ws.Range("A18").Value = "Jodine Personal Savings"
$ws.Range("C18").Value = "Savings"
$ws.Range("D18").Value = "Jodine"

$ws.Range("A19").Value = "Jodine Cash Savings"
$ws.Range("C19").Value = "Savings"
$ws.Range("D19").Value = "Jodine"

# Match the style used in column A for existing account rows (A2:A15)
$ws.Range("A15").Copy()
$ws.Range("A16:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to match the final state
$ws.Range("B14").Select()
